# Adds "Critical temperature", "Critical pressure" and "Critical molar
# volume" rows (31-33) to the COMPOUNDS sheet, one row per property, with
# values for each of the six compounds (columns C-H), the units + source
# columns (L, M) and the parameter name / short name columns (A, B).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COMPOUNDS")

# ---- helpers -------------------------------------------------------------

# Plain text cell (t="s", no cell style) - used for the numeric-looking
# data columns (C:H) which must be stored as text, not as numbers.
function Set-PlainText($cellRef, [string]$value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

# Text cell that also needs the scientific-notation number format (style
# reused from the existing "0.00E+00" style already present in the file).
function Set-PlainTextScientific($cellRef, [string]$value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
    $rng.NumberFormat = "0.00E+00"
}

# Text cell that should carry the same cell style as an existing
# "template" cell (copies format only, value is set beforehand).
function Set-StyledText($cellRef, [string]$value, $templateRef) {
    $rng = $ws.Range($cellRef)
    $rng.Value = $value
    $ws.Range($templateRef).Copy()
    $rng.PasteSpecial(-4122) | Out-Null
}

# ---- row 31: Critical temperature ---------------------------------------

$ws.Range("A31").Value = "Critical temperature"
Set-StyledText "B31" "Tc" "B30"

Set-PlainText "C31" "647.10"
Set-PlainText "D31" "126.19"
Set-PlainText "E31" "190.56"
Set-PlainText "F31" "405.56"
Set-PlainText "G31" "32.94"
Set-PlainText "H31" "456.66"

Set-StyledText "L31" "K" "L30"
Set-StyledText "M31" "CRC, Lange(HCN)" "M30"

# ---- row 32: Critical pressure ------------------------------------------

$ws.Range("A32").Value = "Critical pressure"
Set-StyledText "B32" "pc" "B30"

Set-PlainText "C32" "22.06e6"
Set-PlainText "D32" "3.3958e6"
Set-PlainText "E32" "4.60e6"
Set-PlainText "F32" "11.357e6"
Set-PlainText "G32" "1.2858e6"
Set-PlainText "H32" "5.3905e6"

Set-StyledText "L32" "Pa" "L30"
Set-StyledText "M32" "CRC, Lange(HCN)" "M30"

# ---- row 33: Critical molar volume --------------------------------------

$ws.Range("A33").Value = "Critical molar volume"
Set-StyledText "B33" "Vmc" "B30"

Set-PlainText "C33" "5.6e-5"
Set-PlainTextScientific "D33" "9.0e-5"
Set-PlainText "E33" "9.9e-5"
Set-PlainText "F33" "6.98e-5"
Set-PlainText "G33" "6.5e-5"
Set-PlainText "H33" "1.39e-4"

Set-StyledText "L33" "m3.mol-1" "L30"
Set-StyledText "M33" "CRC, Lange(HCN)" "M30"

# ---- view state -----------------------------------------------------------
# Mirrors the author re-selecting H34 after entering the new rows.
$ws.Range("H34").Select()
